$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 432.4
$ws.Range("I33").Value = 432.4
$ws.Range("K33").Value = 432.4
$ws.Range("M33").Value = -203.4

$ws.Range("H43").Value = 10840.667
$ws.Range("I43").Value = 7276.25
$ws.Range("K43").Value = 7276.25
$ws.Range("M43").Value = -7207.25

$ws.Range("H137").Value = 4819.8
$ws.Range("I137").Value = 2066.6667
$ws.Range("K137").Value = 6200.000100000001
$ws.Range("M137").Value = -3650.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1755.6154
$ws.Range("I2").Value = 1626.6945
$ws.Range("J2").Value = 3302.6667
$ws.Range("K2").Value = 1626.6945
$ws.Range("L2").Value = 3302.6667
$ws.Range("M2").Value = -1513.6945
$ws.Range("N2").Value = -3528.6667

$ws.Range("H4").Value = 182477.45
$ws.Range("I4").Value = 200375.2
$ws.Range("J4").Value = 3500
$ws.Range("K4").Value = 200375.2
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = -200259.2
$ws.Range("N4").Value = -3732

$ws.Range("H5").Value = 514.3125
$ws.Range("I5").Value = 485.75
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 485.75
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = -373.75
$ws.Range("N5").Value = -824

$ws.Range("H32").Value = 6000.4385
$ws.Range("I32").Value = 5226.9243
$ws.Range("K32").Value = 5226.9243
$ws.Range("M32").Value = -4939.9243

$ws.Range("H61").Value = 2252.682
$ws.Range("I61").Value = 2150.4736
$ws.Range("K61").Value = 2150.4736
$ws.Range("M61").Value = -1938.4736

$ws.Range("H63").Value = 3026.2104
$ws.Range("I63").Value = 1919.8572
$ws.Range("K63").Value = 1919.8572
$ws.Range("M63").Value = -1233.8572

$ws.Range("H66").Value = 3026.2104
$ws.Range("I66").Value = 1919.8572
$ws.Range("K66").Value = 9599.286
$ws.Range("M66").Value = -6167.286

$ws.Range("H116").Value = 1755.6154
$ws.Range("I116").Value = 1626.6945
$ws.Range("J116").Value = 3302.6667
$ws.Range("K116").Value = 1626.6945
$ws.Range("L116").Value = 3302.6667
$ws.Range("M116").Value = 667.3054999999999
$ws.Range("N116").Value = -7890.6667

$ws.Range("H132").Value = 5706.4
$ws.Range("I132").Value = 6165.561
$ws.Range("K132").Value = 18496.683
$ws.Range("M132").Value = -15966.683

$ws.Range("H136").Value = 2252.682
$ws.Range("I136").Value = 2150.4736
$ws.Range("K136").Value = 6451.4208
$ws.Range("M136").Value = -3901.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1755.6154
$ws.Range("I3").Value = 1626.6945
$ws.Range("J3").Value = 3302.6667
$ws.Range("K3").Value = 1626.6945
$ws.Range("L3").Value = 3302.6667
$ws.Range("M3").Value = -1512.6945
$ws.Range("N3").Value = -3530.6667

$ws.Range("H4").Value = 514.3125
$ws.Range("I4").Value = 485.75
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 485.75
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = -370.75
$ws.Range("N4").Value = -830

$ws.Range("H94").Value = 915.7368
$ws.Range("I94").Value = 904.7646999999999
$ws.Range("J94").Value = 1009
$ws.Range("K94").Value = 904.7646999999999
$ws.Range("L94").Value = 1009
$ws.Range("M94").Value = -453.7646999999999
$ws.Range("N94").Value = -1911

$ws.Range("H105").Value = 2398.75
$ws.Range("I105").Value = 2200
$ws.Range("J105").Value = 2995
$ws.Range("K105").Value = 2200
$ws.Range("L105").Value = 2995
$ws.Range("M105").Value = -453
$ws.Range("N105").Value = -6489

$ws.Range("H107").Value = 202740
$ws.Range("I107").Value = 251550
$ws.Range("K107").Value = 251550
$ws.Range("M107").Value = -249630

$ws.Range("H134").Value = 301664.7
$ws.Range("I134").Value = 301664.7
$ws.Range("K134").Value = 904994.1000000001
$ws.Range("M134").Value = -902459.1000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1776
$ws.Range("I16").Value = 1532.3
$ws.Range("K16").Value = 1532.3
$ws.Range("M16").Value = -1245.3

$ws.Range("H105").Value = 1772.2778
$ws.Range("I105").Value = 1622.3572
$ws.Range("K105").Value = 1622.3572
$ws.Range("M105").Value = 124.6428000000001

$ws.Range("H107").Value = 1431.05
$ws.Range("I107").Value = 1091.5555
$ws.Range("K107").Value = 1091.5555
$ws.Range("M107").Value = 828.4445000000001

$ws.Range("H113").Value = 1776
$ws.Range("I113").Value = 1532.3
$ws.Range("K113").Value = 1532.3
$ws.Range("M113").Value = 637.7

$ws.Range("H134").Value = 2337.3333
$ws.Range("I134").Value = 1836.1143
$ws.Range("J134").Value = 4843.4287
$ws.Range("K134").Value = 5508.3429
$ws.Range("L134").Value = 14530.2861
$ws.Range("M134").Value = -2973.3429
$ws.Range("N134").Value = -19600.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 217.6
$ws.Range("I12").Value = 31
$ws.Range("J12").Value = 238.33333
$ws.Range("K12").Value = 93
$ws.Range("L12").Value = 714.99999
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = -1060.99999

$ws.Range("H34").Value = 268.94116
$ws.Range("I34").Value = 206.13333
$ws.Range("K34").Value = 618.39999
$ws.Range("M34").Value = -534.39999

$ws.Range("H44").Value = 154.5
$ws.Range("I44").Value = 74.5
$ws.Range("J44").Value = 314.5
$ws.Range("K44").Value = 223.5
$ws.Range("L44").Value = 943.5
$ws.Range("M44").Value = 174.5
$ws.Range("N44").Value = -1739.5

$ws.Range("H55").Value = 7879436.5
$ws.Range("J55").Value = 12505319
$ws.Range("L55").Value = 37515957
$ws.Range("N55").Value = -37516311

$ws.Range("H75").Value = 3103.9
$ws.Range("J75").Value = 4756.6665
$ws.Range("L75").Value = 14269.9995
$ws.Range("N75").Value = -16265.9995

$ws.Range("H78").Value = 3103.9
$ws.Range("J78").Value = 4756.6665
$ws.Range("L78").Value = 42809.9985
$ws.Range("N78").Value = -52793.9985

$ws.Range("H87").Value = 6335.6665
$ws.Range("I87").Value = 6335.6665
$ws.Range("K87").Value = 19006.9995
$ws.Range("M87").Value = -17758.9995

$ws.Range("H90").Value = 6335.6665
$ws.Range("I90").Value = 6335.6665
$ws.Range("K90").Value = 57020.9985
$ws.Range("M90").Value = -50780.9985

$ws.Range("H92").Value = 2092.3333
$ws.Range("J92").Value = 2549.4
$ws.Range("L92").Value = 7648.200000000001
$ws.Range("N92").Value = -10144.2

$ws.Range("H138").Value = 2617.875
$ws.Range("I138").Value = 1998
$ws.Range("J138").Value = 4477.5
$ws.Range("K138").Value = 5994
$ws.Range("L138").Value = 13432.5
$ws.Range("M138").Value = -854
$ws.Range("N138").Value = -23712.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2294.7222
$ws.Range("I102").Value = 2312.0588
$ws.Range("K102").Value = 2312.0588
$ws.Range("M102").Value = -690.0587999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 33264.477
$ws.Range("I40").Value = 7704.9
$ws.Range("J40").Value = 58824.05
$ws.Range("K40").Value = 7704.9
$ws.Range("L40").Value = 58824.05
$ws.Range("M40").Value = -7568.9
$ws.Range("N40").Value = -59096.05

$ws.Range("H55").Value = 91571.82000000001
$ws.Range("I55").Value = 100679
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 100679
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -100506
$ws.Range("N55").Value = -846

$ws.Range("H122").Value = 109204.055
$ws.Range("I122").Value = 136369.73
$ws.Range("K122").Value = 409109.1900000001
$ws.Range("M122").Value = -406659.1900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 36015.75
$ws.Range("I45").Value = 33354.332
$ws.Range("J45").Value = 44000
$ws.Range("K45").Value = 33354.332
$ws.Range("L45").Value = 44000
$ws.Range("M45").Value = -32863.332
$ws.Range("N45").Value = -44982

$ws.Range("H113").Value = 335.35715
$ws.Range("I113").Value = 294.5
$ws.Range("J113").Value = 389.83334
$ws.Range("K113").Value = 883.5
$ws.Range("L113").Value = 1169.50002
$ws.Range("M113").Value = 1286.5
$ws.Range("N113").Value = -5509.500019999999

$ws.Range("H122").Value = 1975.0476
$ws.Range("I122").Value = 1933.2667
$ws.Range("K122").Value = 5799.800099999999
$ws.Range("M122").Value = -3349.800099999999

$ws.Range("H136").Value = 978.2222
$ws.Range("I136").Value = 977.38464
$ws.Range("K136").Value = 2932.15392
$ws.Range("M136").Value = -382.1539199999997
